$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (estimate), C (SE), D (df), E (t.ratio), F (p.value, text)
$rows = @(
    @{ Row = 2;  B = -6.259;             C = 0.234; D = 2007; E = -26.693; F = "1.3605e-134" },
    @{ Row = 3;  B = -4.434;             C = 0.412; D = 2007; E = -10.767; F = "2.5525e-26" },
    @{ Row = 4;  B = -9.457000000000001; C = 0.528; D = 2007; E = -17.925; F = "9.1294e-67" },
    @{ Row = 5;  B = -11.9;              C = 0.574; D = 2007; E = -20.736; F = "1.0636e-86" },
    @{ Row = 6;  B = 0.499;              C = 0.203; D = 2007; E = 2.455;   F = "1.4186e-02" },
    @{ Row = 7;  B = -1.057;             C = 0.105; D = 2007; E = -10.06;  F = "2.9099e-23" },
    @{ Row = 8;  B = -3.664;             C = 0.171; D = 2007; E = -21.401; F = "1.0765e-91" },
    @{ Row = 9;  B = -8.003;             C = 0.361; D = 2007; E = -22.149; F = "1.9838e-97" },
    @{ Row = 10; B = -5.54;              C = 0.414; D = 2007; E = -13.393; F = "3.0552e-39" },
    @{ Row = 11; B = -14.668;            C = 0.503; D = 2007; E = -29.154; F = "4.1901e-156" },
    @{ Row = 12; B = -14.388;            C = 0.53;  D = 2007; E = -27.13;  F = "2.3722e-138" },
    @{ Row = 13; B = -9.563000000000001; C = 0.262; D = 2007; E = -36.539; F = "1.6050e-224" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E

    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Value = "'" + $item.F
}
